$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.291.29"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "2.301.12"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'315.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").Value = "'104.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E7").Value = "  +0.68%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("D10").Value = "'39.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.99%  "

$ws.Range("D11").Value = "'0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").Value = "'8.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").Value = "'0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").Value = "'0.964"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "'15.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("D16").Value = "2.650.40"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "2.298.42"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "42.381.03"
$ws.Range("E18").Value = "  +0.92%  "

$ws.Range("E19").Value = "  -2.75%  "

$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").Value = "'73.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.48%  "

$ws.Range("D22").Value = "'276.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.68%  "

$ws.Range("D23").Value = "'3.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").Value = "'11.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.63%  "

$ws.Range("E25").Value = "  -1.05%  "

$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").Value = "'10.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("E28").Value = "  +3.26%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'36.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'165.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("D32").Value = "'0.0871"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.80%  "

$ws.Range("D33").Value = "'5.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("E34").Value = "  +4.03%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.88%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("E37").Value = "  +3.98%  "

$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("E39").Value = "  +3.48%  "

$ws.Range("E40").Value = "  -0.76%  "

$ws.Range("D41").Value = "'1.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.62%  "

$ws.Range("D42").Value = "'69.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.50%  "

$ws.Range("D43").Value = "'95.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.21%  "

$ws.Range("D44").Value = "'0.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("E45").Value = "  +0.25%  "

$ws.Range("D46").Value = "'81.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.80%  "

$ws.Range("D47").Value = "'12.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.95%  "

$ws.Range("D48").Value = "'113.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D50").Value = "'5.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "

$ws.Range("D51").Value = "1.590.45"
$ws.Range("E51").Value = "  +1.59%  "
